$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:E2")
$rng.NumberFormat = "@"
$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "28.906.91"
$ws.Range("E2").Value = "  -1.85%  "

$rng = $ws.Range("B3:E3")
$rng.NumberFormat = "@"
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.827.03"
$ws.Range("E3").Value = "  -2.34%  "

$rng = $ws.Range("B4:E4")
$rng.NumberFormat = "@"
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.01%  "

$rng = $ws.Range("B5:E5")
$rng.NumberFormat = "@"
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "239.50"
$ws.Range("E5").Value = "  -1.81%  "

$rng = $ws.Range("B6:E6")
$rng.NumberFormat = "@"
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "0.6856"
$ws.Range("E6").Value = "  -2.94%  "

$rng = $ws.Range("B7:E7")
$rng.NumberFormat = "@"
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.01%  "

$rng = $ws.Range("B8:E8")
$rng.NumberFormat = "@"
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "0.07610"
$ws.Range("E8").Value = "  -3.46%  "

$rng = $ws.Range("B9:E9")
$rng.NumberFormat = "@"
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.3013"
$ws.Range("E9").Value = "  -4.49%  "

$rng = $ws.Range("B10:E10")
$rng.NumberFormat = "@"
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "23.42"
$ws.Range("E10").Value = "  -5.08%  "

$rng = $ws.Range("B11:E11")
$rng.NumberFormat = "@"
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07740"
$ws.Range("E11").Value = "  -3.41%  "

$rng = $ws.Range("B12:E12")
$rng.NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.830.65"
$ws.Range("E12").Value = "  -3.45%  "

$rng = $ws.Range("B13:E13")
$rng.NumberFormat = "@"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.038"
$ws.Range("E13").Value = "  -3.46%  "

$rng = $ws.Range("B14:E14")
$rng.NumberFormat = "@"
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "90.14"
$ws.Range("E14").Value = "  -4.37%  "

$rng = $ws.Range("B15:E15")
$rng.NumberFormat = "@"
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.6713"
$ws.Range("E15").Value = "  -4.88%  "

$rng = $ws.Range("B16:E16")
$rng.NumberFormat = "@"
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").Value = "6.415"
$ws.Range("E16").Value = "  -1.20%  "

$rng = $ws.Range("B17:E17")
$rng.NumberFormat = "@"
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.000008263"
$ws.Range("E17").Value = "  -1.31%  "

$rng = $ws.Range("B18:E18")
$rng.NumberFormat = "@"
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "28.911.15"
$ws.Range("E18").Value = "  -2.01%  "

$rng = $ws.Range("B19:E19")
$rng.NumberFormat = "@"
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "242.68"
$ws.Range("E19").Value = "  -5.54%  "

$rng = $ws.Range("B20:E20")
$rng.NumberFormat = "@"
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.094.06"
$ws.Range("E20").Value = "  -2.18%  "

$rng = $ws.Range("B21:E21")
$rng.NumberFormat = "@"
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "12.61"
$ws.Range("E21").Value = "  -4.58%  "

$rng = $ws.Range("B22:E22")
$rng.NumberFormat = "@"
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.08%  "

$rng = $ws.Range("B23:E23")
$rng.NumberFormat = "@"
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "7.387"
$ws.Range("E23").Value = "  -3.26%  "

$rng = $ws.Range("B24:E24")
$rng.NumberFormat = "@"
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.03%  "

$rng = $ws.Range("B25:E25")
$rng.NumberFormat = "@"
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").Value = "0.1465"
$ws.Range("E25").Value = "  -5.90%  "

$rng = $ws.Range("B26:E26")
$rng.NumberFormat = "@"
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "160.82"
$ws.Range("E26").Value = "  -0.14%  "

$rng = $ws.Range("B27:E27")
$rng.NumberFormat = "@"
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "8.704"
$ws.Range("E27").Value = "  -4.09%  "

$rng = $ws.Range("B28:E28")
$rng.NumberFormat = "@"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "18.11"
$ws.Range("E28").Value = "  -3.86%  "

$rng = $ws.Range("B29:E29")
$rng.NumberFormat = "@"
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "1.530"
$ws.Range("E29").Value = "  +1.84%  "

$rng = $ws.Range("B30:E30")
$rng.NumberFormat = "@"
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "4.194"
$ws.Range("E30").Value = "  -3.40%  "

$rng = $ws.Range("B31:E31")
$rng.NumberFormat = "@"
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "4.150"
$ws.Range("E31").Value = "  -2.46%  "

$rng = $ws.Range("B32:E32")
$rng.NumberFormat = "@"
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "1.190"
$ws.Range("E32").Value = "  -1.58%  "

$rng = $ws.Range("B33:E33")
$rng.NumberFormat = "@"
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.05102"
$ws.Range("E33").Value = "  -4.18%  "

$rng = $ws.Range("B34:E34")
$rng.NumberFormat = "@"
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7533"
$ws.Range("E34").Value = "  +0.73%  "

$rng = $ws.Range("B35:E35")
$rng.NumberFormat = "@"
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "1.811"
$ws.Range("E35").Value = "  -4.62%  "

$rng = $ws.Range("B36:E36")
$rng.NumberFormat = "@"
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.140"
$ws.Range("E36").Value = "  -2.82%  "

$rng = $ws.Range("B37:E37")
$rng.NumberFormat = "@"
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.693"
$ws.Range("E37").Value = "  -0.82%  "

$rng = $ws.Range("B38:E38")
$rng.NumberFormat = "@"
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01829"
$ws.Range("E38").Value = "  -2.74%  "

$rng = $ws.Range("B39:E39")
$rng.NumberFormat = "@"
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.218.12"
$ws.Range("E39").Value = "  -3.62%  "

$rng = $ws.Range("B40:E40")
$rng.NumberFormat = "@"
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.682"
$ws.Range("E40").Value = "  -2.59%  "

$rng = $ws.Range("B41:E41")
$rng.NumberFormat = "@"
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.9093"
$ws.Range("E41").Value = "  +0.98%  "

$rng = $ws.Range("B42:E42")
$rng.NumberFormat = "@"
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "108.73"
$ws.Range("E42").Value = "  -0.10%  "

$rng = $ws.Range("B43:E43")
$rng.NumberFormat = "@"
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "0.9996"
$ws.Range("E43").Value = "  -0.05%  "

$rng = $ws.Range("B44:E44")
$rng.NumberFormat = "@"
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.994.37"
$ws.Range("E44").Value = "  -2.18%  "

$rng = $ws.Range("B45:E45")
$rng.NumberFormat = "@"
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.5166"
$ws.Range("E45").Value = "  -0.57%  "

$rng = $ws.Range("B46:E46")
$rng.NumberFormat = "@"
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000122"
$ws.Range("E46").Value = "  -6.10%  "

$rng = $ws.Range("B47:E47")
$rng.NumberFormat = "@"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "9.465"
$ws.Range("E47").Value = "  -0.61%  "

$rng = $ws.Range("B48:E48")
$rng.NumberFormat = "@"
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "5.334"
$ws.Range("E48").Value = "  -10.31%  "

$rng = $ws.Range("B49:E49")
$rng.NumberFormat = "@"
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "62.71"
$ws.Range("E49").Value = "  -12.68%  "

$rng = $ws.Range("B50:E50")
$rng.NumberFormat = "@"
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.721"
$ws.Range("E50").Value = "  -4.98%  "

$rng = $ws.Range("B51:E51")
$rng.NumberFormat = "@"
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "0.4161"
$ws.Range("E51").Value = "  -3.82%  "
